# cryptos.xlsx refresh -- Sun Jun  9 05:55:54 UTC 2024 (GitHub Actions)
# Updates the Price (D) and Volume(1h) (E) columns of the crypto table to the
# latest snapshot. Two coins (EthereumClassic / Binance-PegBSC-USD) also swapped
# rank order, so rows 32 and 33 get their Coin/Link/Price/Volume cells rewritten
# in full rather than patched in place.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 32 & 33: EthereumClassic moved above Binance-PegBSC-USD.
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "26.98"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.72%  "

$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.999"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.08%  "

# All other rows: Price / Volume(1h) refresh only.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "69.321.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.681.11"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.02%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "677.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "158.49"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.02%  "
$ws.Range("E8").Value = "  -1.13%  "
$ws.Range("E9").Value = "  -1.30%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.92"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.435"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.17%  "
$ws.Range("E12").Value = "  -3.38%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.304.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "32.38"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -3.81%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.682.79"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "69.285.41"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.07%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.115"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.57%  "
$ws.Range("E18").Value = "  -2.22%  "
$ws.Range("E19").Value = "  -3.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "467.91"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -3.38%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("E22").Value = "  -2.48%  "
$ws.Range("E23").Value = "  -0.67%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.826.65"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.01%  "
$ws.Range("E25").Value = "  -0.02%  "
$ws.Range("E26").Value = "  -5.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.92"
$ws.Range("D27").Style = "Normal"
$ws.Range("E28").Value = "  -4.62%  "
$ws.Range("E29").Value = "  -1.25%  "
$ws.Range("E30").Value = "  -4.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.62"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -3.24%  "
$ws.Range("E34").Value = "  -4.43%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.674.00"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.88%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.159"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -4.70%  "
$ws.Range("E37").Value = "  -3.21%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.19"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.03%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  -0.01%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "2.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -4.80%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0902"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.30%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "172.05"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.35%  "
$ws.Range("E44").Value = "  -0.89%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "47.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -1.38%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.000280"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.43%  "
$ws.Range("E47").Value = "  -4.74%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "27.96"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -6.00%  "
$ws.Range("E49").Value = "  -4.72%  "
$ws.Range("E50").Value = "  -3.57%  "
$ws.Range("E51").Value = "  -3.09%  "
